# Re-order the data rows (2..26) of the active sheet: each destination row
# ends up holding exactly the values that used to live in a different row
# (full-row relocation, not a per-cell edit). Columns A..AY cover every
# populated cell in this sheet, so snapshotting/restoring that block is
# enough to reproduce the change losslessly (including cells such as J/L
# that are only present on some rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 26
$firstCol = "A"
$lastCol = "AY"

# destination row -> source row (source row's OLD contents become the new
# contents of the destination row)
$mapping = @{
    2  = 13
    3  = 14
    4  = 15
    5  = 2
    6  = 3
    7  = 16
    8  = 4
    9  = 17
    10 = 5
    11 = 18
    12 = 19
    13 = 6
    14 = 7
    15 = 8
    16 = 20
    17 = 21
    18 = 9
    19 = 22
    20 = 23
    21 = 24
    22 = 10
    23 = 25
    24 = 11
    25 = 26
    26 = 12
}

# 1) Snapshot every source row's current values before any writes happen,
#    so later writes never clobber data we still need to read.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range("$firstCol$r`:$lastCol$r")
    $snapshots[$r] = $rng.Value2
}

# 2) Write each destination row from the matching snapshot.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $destRng = $ws.Range("$firstCol$destRow`:$lastCol$destRow")
    $destRng.Value2 = $snapshots[$srcRow]
}
